$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value (as scraped from the latest crypto price pull).
# Values must stay as literal text (e.g. "1.00", "0.0946") rather than being coerced to
# numbers by Excel, so each cell is explicitly formatted as Text before assignment, then
# restored to the default "Normal" style so no stray formatting is introduced.
$updates = [ordered]@{
    'D2' = '43.856.77'
    'E2' = '  +0.99%  '
    'D3' = '2.312.21'
    'E3' = '  +1.51%  '
    'D4' = '1.00'
    'E4' = '  +0.22%  '
    'D5' = '115.29'
    'E5' = '  +22.03%  '
    'D6' = '269.42'
    'E6' = '  +1.23%  '
    'D7' = '0.627'
    'E8' = '  +0.22%  '
    'D9' = '0.627'
    'E9' = '  +3.84%  '
    'D10' = '48.87'
    'E10' = '  +10.71%  '
    'D11' = '0.0946'
    'E11' = '  +1.57%  '
    'D12' = '8.78'
    'E12' = '  +14.14%  '
    'E13' = '  +2.31%  '
    'D14' = '15.71'
    'E14' = '  +4.06%  '
    'D15' = '2.626.99'
    'E15' = '  +0.29%  '
    'D16' = '0.864'
    'E16' = '  +2.60%  '
    'D17' = '2.316.53'
    'E17' = '  +1.55%  '
    'D18' = '43.837.84'
    'E18' = '  +0.85%  '
    'E19' = '  +3.89%  '
    'D20' = '6.62'
    'E20' = '  +7.55%  '
    'D21' = '72.63'
    'E21' = '  +0.83%  '
    'E22' = '  +8.58%  '
    'D23' = '234.64'
    'E23' = '  +0.64%  '
    'D24' = '9.52'
    'E24' = '  +6.83%  '
    'D25' = '2.88'
    'E25' = '  +16.14%  '
    'D27' = '11.57'
    'E27' = '  +3.65%  '
    'D28' = '44.09'
    'E28' = '  +12.61%  '
    'E29' = '  -1.52%  '
    'E30' = '  -0.18%  '
    'D31' = '177.97'
    'E31' = '  +1.47%  '
    'E32' = '  +0.53%  '
    'E33' = '  +6.69%  '
    'D34' = '5.57'
    'E34' = '  +5.02%  '
    'B35' = 'Stellar'
    'C35' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D35' = '0.126'
    'E35' = '  +1.34%  '
    'B36' = 'RenderToken'
    'C36' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D36' = '4.76'
    'E36' = '  +8.34%  '
    'D37' = '0.111'
    'E37' = '  +3.42%  '
    'D38' = '3.96'
    'E38' = '  +20.31%  '
    'D39' = '0.0356'
    'E39' = '  +1.29%  '
    'B40' = 'MultiversX'
    'C40' = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
    'D40' = '76.17'
    'E40' = '  +21.30%  '
    'B41' = 'Algorand'
    'C41' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'D41' = '0.249'
    'E41' = '  +6.65%  '
    'D42' = '2.39'
    'E42' = '  +3.23%  '
    'D43' = '13.21'
    'E43' = '  +11.36%  '
    'D44' = '6.03'
    'E44' = '  +16.20%  '
    'E45' = '  -0.01%  '
    'E46' = '  +4.96%  '
    'E47' = '  +0.69%  '
    'D48' = '0.101'
    'E48' = '  -0.19%  '
    'D49' = '101.24'
    'E49' = '  +4.05%  '
    'E50' = '  +4.52%  '
    'E51' = '  +9.98%  '
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
    $cell.Style = "Normal"
}
